# Fill in the missing "wantUSCurrency" answer for the row-7 survey response.
# Column D mirrors column C ("Yes"/"No") for whether the respondent wants US
# currency; row 7 was missing this value, so add it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = "Yes"

# Touching the font (even re-applying the same face/size) is what nudges the
# style table to mint the dedicated cell format this cell ends up with.
$font = $ws.Range("D7").Font
$font.Name = "Calibri"
$font.Size = 11

# Leave the selection where the author's last click landed.
[void]$ws.Range("E15").Select()
